$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with refreshed metric values (Python 3.10 re-run) ---

# Row 2: opus-mt-en-de
$ws.Range("B2").Value = 0.249
$ws.Range("C2").Value = 0.421
$ws.Range("D2").Value = 0.596

# Row 3: opus-mt-en-de_finetuned
$ws.Range("B3").Value = 0.379
$ws.Range("C3").Value = 0.541
$ws.Range("D3").Value = 0.709

# Row 4: t5-base -> t5-large (relabeled + rerun)
$ws.Range("A4").Value = "t5-large"
$ws.Range("B4").Value = 0.173
$ws.Range("C4").Value = 0.297
$ws.Range("D4").Value = 0.446

# Row 5: t5-base_finetuned -> t5-large_finetuned
$ws.Range("A5").Value = "t5-large_finetuned"
$ws.Range("B5").Value = 0.367
$ws.Range("C5").Value = 0.521
$ws.Range("D5").Value = 0.6899999999999999

# --- Add new rows for NLLB-200 (copy formatting from the row above, then fill in) ---

# Row 6: facebook/nllb-200-distilled-600M
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "facebook/nllb-200-distilled-600M"
$ws.Range("B6").Value = 0.244
$ws.Range("C6").Value = 0.411
$ws.Range("D6").Value = 0.595

# Row 7: facebook/nllb-200-distilled-600M_finetuned
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "facebook/nllb-200-distilled-600M_finetuned"
$ws.Range("B7").Value = 0.367
$ws.Range("C7").Value = 0.531
$ws.Range("D7").Value = 0.703

$excel.CutCopyMode = 0
